$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F2").Value = 60
$ws1.Range("F3").Value = 35
$ws1.Range("F4").Value = 5060
$ws1.Range("F5").Value = 5060
$ws1.Range("F6").Value = 118
$ws1.Range("F7").Value = 147
$ws1.Range("F8").Value = 195
$ws1.Range("F11").Value = 166
$ws1.Range("F12").Value = 8346
$ws1.Range("F16").Value = 607
$ws1.Range("F17").Value = 2518
$ws1.Range("F20").Value = 2297
$ws1.Range("F23").Value = 2519
$ws1.Range("F25").Value = 14
$ws1.Range("F26").Value = 6390
$ws1.Range("F27").Value = 181
$ws1.Range("F29").Value = 133
$ws1.Range("F31").Value = 456
$ws1.Range("F32").Value = 6870
$ws1.Range("F35").Value = 226
$ws1.Range("F36").Value = 12
$ws1.Range("F42").Value = 2521
$ws1.Range("F45").Value = 1126
$ws1.Range("F47").Value = 511
$ws1.Range("F48").Value = 2216
$ws1.Range("F49").Value = 70
$ws2.Range("F2").Value = 9
$ws2.Range("F5").Value = 54
$ws2.Range("F6").Value = 13
$ws2.Range("F12").Value = 154
$ws2.Range("F13").Value = 8
$ws4.Range("F2").Value = 35
$ws4.Range("F3").Value = 5060
$ws4.Range("F4").Value = 5060
$ws4.Range("F5").Value = 118
$ws4.Range("F6").Value = 147
$ws4.Range("F7").Value = 195
$ws4.Range("F10").Value = 166
$ws4.Range("F11").Value = 8346
$ws4.Range("F12").Value = 8346
$ws4.Range("F15").Value = 607
$ws4.Range("F16").Value = 2518
$ws4.Range("F17").Value = 9
$ws4.Range("F20").Value = 2297
$ws4.Range("F21").Value = 54
$ws4.Range("F23").Value = 2519
$ws4.Range("F26").Value = 14
$ws4.Range("F27").Value = 6390
$ws4.Range("F28").Value = 181
$ws4.Range("F31").Value = 133
$ws4.Range("F33").Value = 456
$ws4.Range("F34").Value = 6870
$ws4.Range("F36").Value = 226
$ws4.Range("F41").Value = 2521
$ws4.Range("F43").Value = 1126
$ws4.Range("F45").Value = 511
$ws4.Range("F46").Value = 154
$ws4.Range("F47").Value = 2216
$ws4.Range("F48").Value = 70
$ws4.Range("F49").Value = 8
